# Weekly CompStat report refresh: new crime data collected.
# Updates report header (volume/week dates) and all weekly crime statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: bump volume/issue number and shift the reporting week ---
$ws.Range("A8").Value = "Volume 32   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/10/2025  Through  3/16/2025"

# --- Weekly crime statistics table (rows 15-33) ---
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 300
$ws.Range("N15").Value = -33.333333333333
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 11.111111111111
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 44
$ws.Range("K16").Value = 6.818181818181
$ws.Range("L16").Value = 62.068965517241
$ws.Range("M16").Value = 62.068965517241
$ws.Range("N16").Value = -82.330827067669
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 233.333333333333
$ws.Range("F17").Value = 30
$ws.Range("H17").Value = 57.894736842105
$ws.Range("I17").Value = 68
$ws.Range("J17").Value = 39
$ws.Range("K17").Value = 74.358974358974
$ws.Range("L17").Value = 83.783783783783
$ws.Range("M17").Value = 142.857142857143
$ws.Range("N17").Value = -19.047619047619
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = -2.439024390243
$ws.Range("L18").Value = 37.931034482758
$ws.Range("M18").Value = -11.111111111111
$ws.Range("N18").Value = -91.011235955056
$ws.Range("C19").Value = 35
$ws.Range("D19").Value = 35
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 152
$ws.Range("G19").Value = 120
$ws.Range("H19").Value = 26.666666666666
$ws.Range("I19").Value = 373
$ws.Range("J19").Value = 362
$ws.Range("K19").Value = 3.038674033149
$ws.Range("L19").Value = 2.754820936639
$ws.Range("M19").Value = 19.935691318328
$ws.Range("N19").Value = -73.242467718794
$ws.Range("C20").Value = 3
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 250
$ws.Range("I20").Value = 14
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = 75
$ws.Range("L20").Value = 16.666666666666
$ws.Range("M20").Value = 75
$ws.Range("N20").Value = -85.858585858585
$ws.Range("C21").Value = 53
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = 10.416666666666
$ws.Range("F21").Value = 227
$ws.Range("G21").Value = 176
$ws.Range("H21").Value = 28.977272727272
$ws.Range("I21").Value = 550
$ws.Range("J21").Value = 497
$ws.Range("K21").Value = 10.663983903420
$ws.Range("L21").Value = 16.525423728813
$ws.Range("M21").Value = 27.314814814814
$ws.Range("N21").Value = -76.138828633405
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = -11.764705882352
$ws.Range("L22").Value = -6.25
$ws.Range("M22").Value = 15.384615384615
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 53
$ws.Range("E24").Value = -11.320754716981
$ws.Range("F24").Value = 191
$ws.Range("G24").Value = 235
$ws.Range("H24").Value = -18.723404255319
$ws.Range("I24").Value = 506
$ws.Range("J24").Value = 598
$ws.Range("K24").Value = -15.384615384615
$ws.Range("L24").Value = -0.784313725490
$ws.Range("M24").Value = 50.148367952522
$ws.Range("C25").Value = 43
$ws.Range("D25").Value = 52
$ws.Range("E25").Value = -17.307692307692
$ws.Range("F25").Value = 176
$ws.Range("G25").Value = 228
$ws.Range("H25").Value = -22.807017543859
$ws.Range("I25").Value = 508
$ws.Range("J25").Value = 588
$ws.Range("K25").Value = -13.605442176870
$ws.Range("L25").Value = -6.959706959706
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -5.882352941176
$ws.Range("F26").Value = 59
$ws.Range("H26").Value = 40.476190476190
$ws.Range("I26").Value = 129
$ws.Range("J26").Value = 128
$ws.Range("K26").Value = 0.78125
$ws.Range("L26").Value = 2.380952380952
$ws.Range("M26").Value = 37.234042553191
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 9
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 80
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 33
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 230
$ws.Range("L28").Value = 135.714285714286
$ws.Range("I33").Value = 3

# --- Cells toggling between numeric data and "not applicable" placeholder text ---
# (copy number format from a same-row cell already carrying the desired style,
#  forcing text-type for the literal "0" / "***.* " placeholders)
$ws.Range("D20").Value = 1
$ws.Range("C20").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = 200
$ws.Range("H20").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("M27").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("M27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("D28").Value = 4
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = -50
$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "0"
$ws.Range("M31").Copy()
$ws.Range("C31").PasteSpecial(-4122)

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("M31").Copy()
$ws.Range("D31").PasteSpecial(-4122)

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("M31").Copy()
$ws.Range("E31").PasteSpecial(-4122)

